$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two newest bank-statement movements arrived; insert two rows at the very
# top of the ledger (rows shift down, the 5-row empty buffer at the bottom
# grows back automatically) and fill them in.
$ws.Rows("1:2").Insert()

# Row 3 (the old row 1) already kept its date style (s=1) on column A after
# the insert; clone that number format onto the two freshly inserted blank
# rows so the new dates render the same way.
$ws.Range("A3").Copy()
$ws.Range("A1:A2").PasteSpecial(-4122) # xlPasteFormats

# --- Row 1 : 2014-02-12 MASTERCARD payment (fill A:F, saldo G comes later) ---
$ws.Range("A1").Value = 41682
$ws.Range("B1").Value = "13118145-MASTERCARD-RA-518114000072"
$ws.Range("C1").Value = "D"
$ws.Range("D1").Value = "0007629440"
$ws.Range("E1").Value = "SERVICIOS CENTRALES"
$ws.Range("F1").Value = "718.86  "

# --- Row 2 : 2014-02-10 DEPOSITO (all columns, including saldo) ---
$ws.Range("A2").Value = 41680
$ws.Range("B2").Value = "DEPOSITO"
$ws.Range("C2").Value = "C"
$ws.Range("D2").Value = "0007829636"
$ws.Range("E2").Value = "TENA"
$ws.Range("F2").Value = "115.00  "
$ws.Range("G2").Value = "2612.81"

# --- back to Row 1 : running balance after the newest movement ---
$ws.Range("G1").Value = "1893.95"

# The PHP-array helper formula only ever lives in H1 (the newest row); the
# row that used to be H1 never had its own H formula, so just drop whatever
# the insert/shift left behind in H3 and (re)build it fresh in H1, with the
# extra 'mo_borrado_logico' fragment from this edit.
$ws.Range("H3").ClearContents()
$ws.Range("H1").Formula = '=CONCATENATE("array(''mo_fecha'' => new \DateTime(''",TEXT(A1,"yyyy-mm-dd"),"''), ''mo_concepto'' => ''",B1,"'', ''mo_tipo'' => ''",C1,"'', ''mo_documento'' => ''",D1,"'', ''mo_oficina'' => ''",E1,"'', ''mo_monto'' => ",TRIM(F1),", ''mo_saldo'' => ",G1,", ''mo_fecha_crea'' => new \DateTime(''",TEXT(NOW(),"yyyy-mm-dd H:m:s"),"''), ''mo_quien_crea'' => 1, ''mo_fecha_modifica'' => NULL, ''mo_quien_modifica'' => NULL ''mo_borrado_logico'' => false),")'

# Restore the usual selection habit of landing on the most recent saldo cell.
$ws.Range("G2").Select()

Write-Output "done"
